$wb = $excel.ActiveWorkbook

$wsDescription = $wb.Worksheets.Item("Description")
$wsBlackbox    = $wb.Worksheets.Item("Blackbox")
$wsWhitebox    = $wb.Worksheets.Item("Whitebox")

# --- Rename Blackbox test IDs: BB_Testvalidate{N} -> BB_validate{N} ---
$wsBlackbox.Range("A3").Value = "BB_validate1"
$wsBlackbox.Range("A4").Value = "BB_validate1"
$wsBlackbox.Range("A5").Value = "BB_validate3"
$wsBlackbox.Range("A6").Value = "BB_validate4"
$wsBlackbox.Range("A7").Value = "BB_validate5"
$wsBlackbox.Range("A8").Value = "BB_validate6"

# --- Rename Whitebox test IDs: WB_Testvalidate{N} -> WB_validate{N} ---
$wsWhitebox.Range("A3").Value = "WB_validate1"
$wsWhitebox.Range("A4").Value = "WB_validate2"
$wsWhitebox.Range("A5").Value = "WB_validate3"

# --- Fix up the "Setup" description text on the Description sheet: a stray
#     linebreak was introduced in the middle of the long run of spaces that
#     separates steps 2 and 3 (content otherwise unchanged). ---
$wsDescription.Range("A12").Value = "1. Recognize the constraints and the requirements of the validate function.                                                                                                                                                                                                     2. Testing both valid and invalid test cases.                                                                                                                            `n 3. Compareing the actual results to the predicted results."

# --- Selection / active-sheet bookkeeping to match the saved view state ---
$wsDescription.Activate()
$wsDescription.Range("L16").Select() | Out-Null

$wsWhitebox.Activate()
$wsWhitebox.Range("A5").Select() | Out-Null

$wsBlackbox.Activate()
$wsBlackbox.Range("A8").Select() | Out-Null
